$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-06-07 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-08 Saturday", 2)

$d.Content.Find.Execute("57×18=1026", $true, $false, $false, $false, $false, $true, 1, $false, "60×94=5640", 2)
$d.Content.Find.Execute("98×39=3822", $true, $false, $false, $false, $false, $true, 1, $false, "69×48=3312", 2)
$d.Content.Find.Execute("36×26=936", $true, $false, $false, $false, $false, $true, 1, $false, "52×13=676", 2)
$d.Content.Find.Execute("71×91=6461", $true, $false, $false, $false, $false, $true, 1, $false, "46×94=4324", 2)
$d.Content.Find.Execute("71×70=4970", $true, $false, $false, $false, $false, $true, 1, $false, "40×41=1640", 2)
$d.Content.Find.Execute("37×43=1591", $true, $false, $false, $false, $false, $true, 1, $false, "66×16=1056", 2)
$d.Content.Find.Execute("94×96=9024", $true, $false, $false, $false, $false, $true, 1, $false, "44×15=660", 2)
$d.Content.Find.Execute("61×13=793", $true, $false, $false, $false, $false, $true, 1, $false, "63×38=2394", 2)
$d.Content.Find.Execute("40×18=720", $true, $false, $false, $false, $false, $true, 1, $false, "17×58=986", 2)
$d.Content.Find.Execute("55×13=715", $true, $false, $false, $false, $false, $true, 1, $false, "73×21=1533", 2)
$d.Content.Find.Execute("28×73=2044", $true, $false, $false, $false, $false, $true, 1, $false, "95×70=6650", 2)
$d.Content.Find.Execute("96×47=4512", $true, $false, $false, $false, $false, $true, 1, $false, "12×90=1080", 2)
$d.Content.Find.Execute("33×69=2277", $true, $false, $false, $false, $false, $true, 1, $false, "95×87=8265", 2)
$d.Content.Find.Execute("54×77=4158", $true, $false, $false, $false, $false, $true, 1, $false, "50×71=3550", 2)
$d.Content.Find.Execute("49×39=1911", $true, $false, $false, $false, $false, $true, 1, $false, "80×64=5120", 2)
$d.Content.Find.Execute("59×13=767", $true, $false, $false, $false, $false, $true, 1, $false, "28×59=1652", 2)
$d.Content.Find.Execute("53×94=4982", $true, $false, $false, $false, $false, $true, 1, $false, "56×38=2128", 2)
$d.Content.Find.Execute("19×70=1330", $true, $false, $false, $false, $false, $true, 1, $false, "41×15=615", 2)
$d.Content.Find.Execute("51×28=1428", $true, $false, $false, $false, $false, $true, 1, $false, "87×38=3306", 2)
$d.Content.Find.Execute("94×60=5640", $true, $false, $false, $false, $false, $true, 1, $false, "36×16=576", 2)
$d.Content.Find.Execute("90×24=2160", $true, $false, $false, $false, $false, $true, 1, $false, "74×38=2812", 2)
$d.Content.Find.Execute("97×15=1455", $true, $false, $false, $false, $false, $true, 1, $false, "83×67=5561", 2)
$d.Content.Find.Execute("55×88=4840", $true, $false, $false, $false, $false, $true, 1, $false, "84×92=7728", 2)
$d.Content.Find.Execute("70×88=6160", $true, $false, $false, $false, $false, $true, 1, $false, "42×53=2226", 2)
$d.Content.Find.Execute("11×77=847", $true, $false, $false, $false, $false, $true, 1, $false, "14×67=938", 2)
